# Auto-generated script applying scheduled-runner updates to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1223.1666
$ws.Cells.Item(19, 9).Value = 1520.4
$ws.Cells.Item(19, 10).Value = 1108.8462
$ws.Cells.Item(19, 11).Value = 1520.4
$ws.Cells.Item(19, 12).Value = 1108.8462
$ws.Cells.Item(19, 13).Value = -1345.4
$ws.Cells.Item(19, 14).Value = -1458.8462
$ws.Cells.Item(33, 8).Value = 703.25
$ws.Cells.Item(33, 9).Value = 120.8
$ws.Cells.Item(33, 11).Value = 120.8
$ws.Cells.Item(33, 13).Value = 108.2
$ws.Cells.Item(82, 8).Value = 3566.6667
$ws.Cells.Item(85, 8).Value = 3566.6667
$ws.Cells.Item(92, 8).Value = 1667.4231
$ws.Cells.Item(92, 9).Value = 1335.15
$ws.Cells.Item(92, 10).Value = 2775
$ws.Cells.Item(92, 11).Value = 1335.15
$ws.Cells.Item(92, 12).Value = 2775
$ws.Cells.Item(92, 13).Value = -87.15000000000009
$ws.Cells.Item(92, 14).Value = -5271
$ws.Cells.Item(113, 8).Value = 54255.21
$ws.Cells.Item(113, 9).Value = 144042.72
$ws.Cells.Item(113, 10).Value = 1879.1666
$ws.Cells.Item(113, 11).Value = 144042.72
$ws.Cells.Item(113, 12).Value = 1879.1666
$ws.Cells.Item(113, 13).Value = -140788.72
$ws.Cells.Item(113, 14).Value = -8387.1666
$ws.Cells.Item(116, 8).Value = 3896.6667
$ws.Cells.Item(116, 9).Value = 2850
$ws.Cells.Item(116, 10).Value = 5990
$ws.Cells.Item(116, 11).Value = 2850
$ws.Cells.Item(116, 12).Value = 5990
$ws.Cells.Item(116, 13).Value = 592
$ws.Cells.Item(116, 14).Value = -12874
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 154.55556
$ws.Cells.Item(5, 9).Value = 141.57143
$ws.Cells.Item(5, 11).Value = 141.57143
$ws.Cells.Item(5, 13).Value = -29.57142999999999
$ws.Cells.Item(32, 8).Value = 29068.62
$ws.Cells.Item(32, 9).Value = 5474.3125
$ws.Cells.Item(32, 11).Value = 5474.3125
$ws.Cells.Item(32, 13).Value = -5187.3125
$ws.Cells.Item(35, 8).Value = 1230
$ws.Cells.Item(35, 9).Value = 1230
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1230
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -824
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(36, 8).Value = 1000
$ws.Cells.Item(36, 9).Value = 1000
$ws.Cells.Item(36, 11).Value = 1000
$ws.Cells.Item(36, 13).Value = -654
$ws.Cells.Item(61, 8).Value = 2206.1304
$ws.Cells.Item(61, 9).Value = 1247.4445
$ws.Cells.Item(61, 11).Value = 1247.4445
$ws.Cells.Item(61, 13).Value = -1035.4445
$ws.Cells.Item(136, 8).Value = 2206.1304
$ws.Cells.Item(136, 9).Value = 1247.4445
$ws.Cells.Item(136, 11).Value = 3742.3335
$ws.Cells.Item(136, 13).Value = -1192.3335
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 154.55556
$ws.Cells.Item(4, 9).Value = 141.57143
$ws.Cells.Item(4, 11).Value = 141.57143
$ws.Cells.Item(4, 13).Value = -26.57142999999999
$ws.Cells.Item(36, 8).Value = 31018.5
$ws.Cells.Item(36, 9).Value = 31018.5
$ws.Cells.Item(36, 11).Value = 31018.5
$ws.Cells.Item(36, 13).Value = -30484.5
$ws.Cells.Item(99, 8).Value = 2046.6666
$ws.Cells.Item(99, 9).Value = 2395
$ws.Cells.Item(99, 10).Value = 1872.5
$ws.Cells.Item(99, 11).Value = 2395
$ws.Cells.Item(99, 12).Value = 1872.5
$ws.Cells.Item(99, 13).Value = -897
$ws.Cells.Item(99, 14).Value = -4868.5
$ws.Cells.Item(134, 8).Value = 1536.95
$ws.Cells.Item(134, 9).Value = 1510.0286
$ws.Cells.Item(134, 11).Value = 4530.085800000001
$ws.Cells.Item(134, 13).Value = -1995.085800000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1030.1
$ws.Cells.Item(16, 9).Value = 951.8333
$ws.Cells.Item(16, 10).Value = 1147.5
$ws.Cells.Item(16, 11).Value = 951.8333
$ws.Cells.Item(16, 12).Value = 1147.5
$ws.Cells.Item(16, 13).Value = -664.8333
$ws.Cells.Item(16, 14).Value = -1721.5
$ws.Cells.Item(31, 8).Value = 1467.5797
$ws.Cells.Item(31, 9).Value = 1014.9524
$ws.Cells.Item(31, 10).Value = 2171.6667
$ws.Cells.Item(31, 11).Value = 1014.9524
$ws.Cells.Item(31, 12).Value = 2171.6667
$ws.Cells.Item(31, 13).Value = -719.9524
$ws.Cells.Item(31, 14).Value = -2761.6667
$ws.Cells.Item(34, 8).Value = 1467.5797
$ws.Cells.Item(34, 9).Value = 1014.9524
$ws.Cells.Item(34, 10).Value = 2171.6667
$ws.Cells.Item(34, 11).Value = 1014.9524
$ws.Cells.Item(34, 12).Value = 2171.6667
$ws.Cells.Item(34, 13).Value = -812.9524
$ws.Cells.Item(34, 14).Value = -2575.6667
$ws.Cells.Item(99, 8).Value = 11499.637
$ws.Cells.Item(99, 9).Value = 2200.5
$ws.Cells.Item(99, 10).Value = 22658.6
$ws.Cells.Item(99, 11).Value = 2200.5
$ws.Cells.Item(99, 12).Value = 22658.6
$ws.Cells.Item(99, 13).Value = -702.5
$ws.Cells.Item(99, 14).Value = -25654.6
$ws.Cells.Item(113, 8).Value = 1030.1
$ws.Cells.Item(113, 9).Value = 951.8333
$ws.Cells.Item(113, 10).Value = 1147.5
$ws.Cells.Item(113, 11).Value = 951.8333
$ws.Cells.Item(113, 12).Value = 1147.5
$ws.Cells.Item(113, 13).Value = 1218.1667
$ws.Cells.Item(113, 14).Value = -5487.5
$ws.Cells.Item(122, 8).Value = 3326.3635
$ws.Cells.Item(122, 9).Value = 3359
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 10077
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -7627
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(126, 8).Value = 11499.637
$ws.Cells.Item(126, 9).Value = 2200.5
$ws.Cells.Item(126, 10).Value = 22658.6
$ws.Cells.Item(126, 11).Value = 6601.5
$ws.Cells.Item(126, 12).Value = 67975.79999999999
$ws.Cells.Item(126, 13).Value = -4131.5
$ws.Cells.Item(126, 14).Value = -72915.79999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 648318.1
$ws.Cells.Item(37, 10).Value = 648318.1
$ws.Cells.Item(37, 12).Value = 1944954.3
$ws.Cells.Item(37, 14).Value = -1945178.3
$ws.Cells.Item(64, 8).Value = 1998.1666
$ws.Cells.Item(67, 8).Value = 1998.1666
$ws.Cells.Item(131, 8).Value = 894.4861
$ws.Cells.Item(131, 10).Value = 938.7077
$ws.Cells.Item(131, 12).Value = 2816.1231
$ws.Cells.Item(131, 14).Value = -12896.1231
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3976
$ws.Cells.Item(80, 9).Value = 5350
$ws.Cells.Item(80, 11).Value = 5350
$ws.Cells.Item(80, 13).Value = -4352
$ws.Cells.Item(83, 8).Value = 3976
$ws.Cells.Item(83, 9).Value = 5350
$ws.Cells.Item(83, 11).Value = 26750
$ws.Cells.Item(83, 13).Value = -21758
$ws.Cells.Item(102, 8).Value = 752829.1
$ws.Cells.Item(102, 9).Value = 7490
$ws.Cells.Item(102, 10).Value = 1001275.5
$ws.Cells.Item(102, 11).Value = 7490
$ws.Cells.Item(102, 12).Value = 1001275.5
$ws.Cells.Item(102, 13).Value = -5868
$ws.Cells.Item(102, 14).Value = -1004519.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2126.2307
$ws.Cells.Item(82, 9).Value = 1542.1111
$ws.Cells.Item(82, 11).Value = 1542.1111
$ws.Cells.Item(82, 13).Value = -1181.1111
$ws.Cells.Item(85, 8).Value = 2126.2307
$ws.Cells.Item(85, 9).Value = 1542.1111
$ws.Cells.Item(85, 11).Value = 1542.1111
$ws.Cells.Item(85, 13).Value = -294.1111000000001
$ws.Cells.Item(87, 8).Value = 37426.668
$ws.Cells.Item(87, 10).Value = 37426.668
$ws.Cells.Item(87, 12).Value = 37426.668
$ws.Cells.Item(87, 14).Value = -39672.668
$ws.Cells.Item(90, 8).Value = 37426.668
$ws.Cells.Item(90, 10).Value = 37426.668
$ws.Cells.Item(90, 12).Value = 112280.004
$ws.Cells.Item(90, 14).Value = -123512.004
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 35997.5
$ws.Cells.Item(57, 10).Value = 35997.5
$ws.Cells.Item(57, 12).Value = 35997.5
$ws.Cells.Item(57, 14).Value = -37505.5
$ws.Cells.Item(122, 8).Value = 1699.8182
$ws.Cells.Item(122, 9).Value = 1500
$ws.Cells.Item(122, 10).Value = 2049.5
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 12).Value = 6148.5
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(122, 14).Value = -11048.5
$ws.Cells.Item(126, 8).Value = 2980.8333
$ws.Cells.Item(126, 9).Value = 3660
$ws.Cells.Item(126, 10).Value = 2301.6667
$ws.Cells.Item(126, 11).Value = 10980
$ws.Cells.Item(126, 12).Value = 6905.000100000001
$ws.Cells.Item(126, 13).Value = -8510
$ws.Cells.Item(126, 14).Value = -11845.0001
$ws.Cells.Item(136, 8).Value = 1432.5
$ws.Cells.Item(136, 9).Value = 662.36365
$ws.Cells.Item(136, 11).Value = 1987.09095
$ws.Cells.Item(136, 13).Value = 562.90905
